$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Access" column for the existing users (Marko Majkic, Nellie Solis,
# Brianna Spendlove, Donna Davis) shares one string -- "boulevard" became "lvb".
# Using Replace() edits the shared text itself so every row referencing it
# updates together, exactly like the source change.
$ws.Cells.Replace("boulevard", "lvb")

# --- Row 6: Jairo.Contreras ---
$ws.Range("A6").Value = "Jairo.Contreras"
$ws.Range("B6").Value = "jairo.contreras@hgv.com"
$ws.Range("C6").Value = "golive2025"
$ws.Range("D6").Value = "cancun"
$ws.Rows.Item(6).RowHeight = 14.9
$ws.Hyperlinks.Add($ws.Range("B6"), "mailto:jairo.contreras@hgv.com", "", "", "jairo.contreras@hgv.com")
$ws.Range("B6").Font.Underline = $false
$ws.Range("B6").Font.Color = 0

# --- Row 7: Miriam Ghiasi ---
$ws.Range("A7").Value = "Miriam Ghiasi"
$ws.Range("B7").Value = "miriam.ghiasi@hgv.com"
$ws.Range("C7").Value = "golive2025"
$ws.Range("D7").Value = "lvb"
$ws.Rows.Item(7).RowHeight = 14.9
$ws.Hyperlinks.Add($ws.Range("B7"), "mailto:miriam.ghiasi@hgv.com", "", "", "miriam.ghiasi@hgv.com")
$ws.Range("B7").Font.Underline = $false
$ws.Range("B7").Font.Color = 0

# --- Row 8: Annie Solseng ---
$ws.Range("A8").Value = "Annie Solseng"
$ws.Range("B8").Value = "annie.solseng@hgv.com"
$ws.Range("C8").Value = "golive2025"
$ws.Range("D8").Value = "elara"
$ws.Rows.Item(8).RowHeight = 14.9
$ws.Hyperlinks.Add($ws.Range("B8"), "mailto:annie.solseng@hgv.com", "", "", "annie.solseng@hgv.com")
$ws.Range("B8").Font.Underline = $false
$ws.Range("B8").Font.Color = 0

# Match the author's final selection
$ws.Range("E10").Select()
